# Actualizacion 11 de Mayo - Manana
# Updates the E..K (Aprobados/Reprobados counts, percentages, and average)
# columns for several rows across the three partial-exam worksheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("1er Parcial")
$ws.Range("E6").Value = 20; $ws.Range("F6").Value = 9; $ws.Range("G6").Value = 68.97; $ws.Range("H6").Value = 31.03; $ws.Range("I6").Value = 8.6; $ws.Range("J6").Value = 9; $ws.Range("K6").Value = 31.03
$ws.Range("E7").Value = 21; $ws.Range("F7").Value = 8; $ws.Range("G7").Value = 72.41; $ws.Range("H7").Value = 27.59; $ws.Range("I7").Value = 8.6; $ws.Range("J7").Value = 8; $ws.Range("K7").Value = 27.59
$ws.Range("E8").Value = 18; $ws.Range("F8").Value = 4; $ws.Range("G8").Value = 81.81999999999999; $ws.Range("H8").Value = 18.18; $ws.Range("I8").Value = 8.300000000000001; $ws.Range("J8").Value = 4; $ws.Range("K8").Value = 18.18
$ws.Range("E10").Value = 39; $ws.Range("F10").Value = 0; $ws.Range("G10").Value = 100; $ws.Range("H10").Value = 0; $ws.Range("I10").Value = 8.199999999999999; $ws.Range("J10").Value = 0; $ws.Range("K10").Value = 0
$ws.Range("E11").Value = 37; $ws.Range("F11").Value = 0; $ws.Range("G11").Value = 100; $ws.Range("H11").Value = 0; $ws.Range("I11").Value = 7.5; $ws.Range("J11").Value = 0; $ws.Range("K11").Value = 0
$ws.Range("E12").Value = 36; $ws.Range("F12").Value = 0; $ws.Range("G12").Value = 100; $ws.Range("H12").Value = 0; $ws.Range("I12").Value = 6.6; $ws.Range("J12").Value = 0; $ws.Range("K12").Value = 0
$ws.Range("E13").Value = 29; $ws.Range("F13").Value = 0; $ws.Range("G13").Value = 100; $ws.Range("H13").Value = 0; $ws.Range("I13").Value = 6.4; $ws.Range("J13").Value = 0; $ws.Range("K13").Value = 0
$ws.Range("E14").Value = 29; $ws.Range("F14").Value = 0; $ws.Range("G14").Value = 100; $ws.Range("H14").Value = 0; $ws.Range("I14").Value = 6.4; $ws.Range("J14").Value = 0; $ws.Range("K14").Value = 0

$ws = $wb.Worksheets.Item("2o Parcial")
$ws.Range("E6").Value = 20; $ws.Range("F6").Value = 9; $ws.Range("G6").Value = 68.97; $ws.Range("H6").Value = 31.03; $ws.Range("I6").Value = 8.4; $ws.Range("J6").Value = 9; $ws.Range("K6").Value = 31.03
$ws.Range("E7").Value = 19; $ws.Range("F7").Value = 10; $ws.Range("G7").Value = 65.52; $ws.Range("H7").Value = 34.48; $ws.Range("I7").Value = 7.6; $ws.Range("J7").Value = 10; $ws.Range("K7").Value = 34.48
$ws.Range("E8").Value = 17; $ws.Range("F8").Value = 5; $ws.Range("G8").Value = 77.27; $ws.Range("H8").Value = 22.73; $ws.Range("I8").Value = 7.4; $ws.Range("J8").Value = 5; $ws.Range("K8").Value = 22.73
$ws.Range("E9").Value = 23; $ws.Range("F9").Value = 16; $ws.Range("G9").Value = 58.97; $ws.Range("H9").Value = 41.03; $ws.Range("I9").Value = 7; $ws.Range("J9").Value = 16; $ws.Range("K9").Value = 41.03
$ws.Range("E10").Value = 39; $ws.Range("F10").Value = 0; $ws.Range("G10").Value = 100; $ws.Range("H10").Value = 0; $ws.Range("I10").Value = 8; $ws.Range("J10").Value = 0; $ws.Range("K10").Value = 0
$ws.Range("E11").Value = 37; $ws.Range("F11").Value = 0; $ws.Range("G11").Value = 100; $ws.Range("H11").Value = 0; $ws.Range("I11").Value = 7.5; $ws.Range("J11").Value = 0; $ws.Range("K11").Value = 0
$ws.Range("E12").Value = 36; $ws.Range("F12").Value = 0; $ws.Range("G12").Value = 100; $ws.Range("H12").Value = 0; $ws.Range("I12").Value = 6.6; $ws.Range("J12").Value = 0; $ws.Range("K12").Value = 0
$ws.Range("E13").Value = 29; $ws.Range("F13").Value = 0; $ws.Range("G13").Value = 100; $ws.Range("H13").Value = 0; $ws.Range("I13").Value = 6.6; $ws.Range("J13").Value = 0; $ws.Range("K13").Value = 0
$ws.Range("E14").Value = 29; $ws.Range("F14").Value = 0; $ws.Range("G14").Value = 100; $ws.Range("H14").Value = 0; $ws.Range("I14").Value = 6.4; $ws.Range("J14").Value = 0; $ws.Range("K14").Value = 0
$ws.Range("E18").Value = 21; $ws.Range("F18").Value = 16; $ws.Range("G18").Value = 56.76; $ws.Range("H18").Value = 43.24; $ws.Range("J18").Value = 16; $ws.Range("K18").Value = 43.24
$ws.Range("E19").Value = 8; $ws.Range("F19").Value = 28; $ws.Range("G19").Value = 22.22; $ws.Range("H19").Value = 77.78; $ws.Range("I19").Value = 8.9; $ws.Range("J19").Value = 28; $ws.Range("K19").Value = 77.78
$ws.Range("E22").Value = 12; $ws.Range("F22").Value = 10; $ws.Range("G22").Value = 54.55; $ws.Range("H22").Value = 45.45; $ws.Range("I22").Value = 8.699999999999999; $ws.Range("J22").Value = 10; $ws.Range("K22").Value = 45.45

$ws = $wb.Worksheets.Item("3er Parcial")
$ws.Range("E6").Value = 20; $ws.Range("F6").Value = 9; $ws.Range("G6").Value = 68.97; $ws.Range("H6").Value = 31.03; $ws.Range("I6").Value = 8.699999999999999; $ws.Range("J6").Value = 9; $ws.Range("K6").Value = 31.03
$ws.Range("E7").Value = 21; $ws.Range("F7").Value = 8; $ws.Range("G7").Value = 72.41; $ws.Range("H7").Value = 27.59; $ws.Range("I7").Value = 8.300000000000001; $ws.Range("J7").Value = 8; $ws.Range("K7").Value = 27.59
$ws.Range("E8").Value = 18; $ws.Range("F8").Value = 4; $ws.Range("G8").Value = 81.81999999999999; $ws.Range("H8").Value = 18.18; $ws.Range("I8").Value = 8; $ws.Range("J8").Value = 4; $ws.Range("K8").Value = 18.18
$ws.Range("E10").Value = 39; $ws.Range("F10").Value = 0; $ws.Range("G10").Value = 100; $ws.Range("H10").Value = 0; $ws.Range("I10").Value = 8.199999999999999; $ws.Range("J10").Value = 0; $ws.Range("K10").Value = 0
$ws.Range("E11").Value = 37; $ws.Range("F11").Value = 0; $ws.Range("G11").Value = 100; $ws.Range("H11").Value = 0; $ws.Range("I11").Value = 7.6; $ws.Range("J11").Value = 0; $ws.Range("K11").Value = 0
$ws.Range("E12").Value = 36; $ws.Range("F12").Value = 0; $ws.Range("G12").Value = 100; $ws.Range("H12").Value = 0; $ws.Range("I12").Value = 6.6; $ws.Range("J12").Value = 0; $ws.Range("K12").Value = 0
$ws.Range("E13").Value = 29; $ws.Range("F13").Value = 0; $ws.Range("G13").Value = 100; $ws.Range("H13").Value = 0; $ws.Range("I13").Value = 6.6; $ws.Range("J13").Value = 0; $ws.Range("K13").Value = 0
$ws.Range("E14").Value = 29; $ws.Range("F14").Value = 0; $ws.Range("G14").Value = 100; $ws.Range("H14").Value = 0; $ws.Range("I14").Value = 6.4; $ws.Range("J14").Value = 0; $ws.Range("K14").Value = 0
